$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad / Changed date) for rows 2 through 13
# from serial date 45207 (2023-10-08) to 45208 (2023-10-09)
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}
